$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 27 for the LaTeX image option
$ws.Rows.Item(27).Insert()
$ws.Cells.Item(27, 1).Value = "Image with LaTeX installed"
$ws.Cells.Item(27, 2).Value = "latexImage"

# 2. Update AWS instance type example value (was subnet id placeholder, now an instance type)
$ws.Cells.Item(33, 3).Value = "m4.xlarge"

# 3. Append new rows for AWS EFS shared storage options
$ws.Cells.Item(36, 1).Value = "AWS EFS/sharedStorage ID"
$ws.Cells.Item(36, 2).Value = "sharedStorage"

$ws.Cells.Item(37, 1).Value = "Mount of Amazon EFS/sharedStorage"
$ws.Cells.Item(37, 2).Value = "sharedStorageMount"
$ws.Cells.Item(37, 3).Value = "/mnt/shared"

# 4. Update selection / scroll position to reflect the new view
$ws.Range("C36").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
